$d = $word.ActiveDocument

# --- Change 1: simplify "mas provavelmente uma mistura de mobile + desktop" ---
# Track revisions while editing so Word does not opportunistically coalesce
# this run back into its identically-formatted neighbours ("definida," and
# the following "."); accepting the tracked change afterwards leaves the
# paragraph's original run layout untouched apart from the edited run.
$d.TrackRevisions = $true

$r1 = $d.Content
$r1.Find.Execute(" mas provavelmente uma mistura de mobile + desktop", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
if ($r1.Find.Found) {
    $r1.Text = " mas provavelmente de mobile"
}

$d.TrackRevisions = $false
$d.Revisions.AcceptAll()

# --- Change 2: collapse "Multiplayer online; 2D; " + "Por" + " turnos." (with
# gramStart/gramEnd proofErr markers around "Por") into a single run reading
# "Multiplayer online; 2D; Por turnos." ---
# The visible text is already exactly the target string, so a direct
# Find/Replace with identical text is a no-op; round-trip through a marker
# value first so Word actually rewrites (and so merges/clears proofErr on)
# the run.
$r2 = $d.Content
$r2.Find.Execute("Multiplayer online; 2D; Por turnos.", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
if ($r2.Find.Found) {
    $r2.Text = "Multiplayer online; 2D; Por turnos.###TMP###"
}

$r3 = $d.Content
$r3.Find.Execute("Multiplayer online; 2D; Por turnos.###TMP###", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
if ($r3.Find.Found) {
    $r3.Text = "Multiplayer online; 2D; Por turnos."
}
